$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 <- old row 16 values
$ws.Range("A7").Value = "Andrew Wiggins"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Golden State Warriors"

# Row 10 <- old row 11 values
$ws.Range("A10").Value = "Jusuf Nurkic"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Phoenix Suns"

# Row 11 <- old row 14 values
$ws.Range("A11").Value = "Ivica Zubac"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "LA Clippers"

# Row 14 <- old row 10 values
$ws.Range("A14").Value = "Draymond Green"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Golden State Warriors"

# Row 16 <- old row 17 values
$ws.Range("A16").Value = "Wendell Carter Jr."
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Orlando Magic"

# Row 17 <- old row 19 values
$ws.Range("A17").Value = "Paul George"
$ws.Range("B17").Value = "SG,SF,PF"
$ws.Range("C17").Value = "Philadelphia 76ers"

# Row 19 <- old row 7 values
$ws.Range("A19").Value = "Giannis Antetokounmpo"
$ws.Range("B19").Value = "PF,C"
$ws.Range("C19").Value = "Milwaukee Bucks"
